$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-13 (A:T) with the refreshed TPM-derived values
# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Leap2"
$ws.Cells.Item(2, 3).Value = "Ghsr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.093697666666667
$ws.Cells.Item(2, 8).Value = 9.281093
$ws.Cells.Item(2, 9).Value = 0.3054078832615814
$ws.Cells.Item(2, 10).Value = 0.3054078832615814
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.5447316666666667
$ws.Cells.Item(2, 14).Value = 1.634195
$ws.Cells.Item(2, 15).Value = 0.1484165462704666
$ws.Cells.Item(2, 16).Value = 0.1484165462704666
$ws.Cells.Item(2, 17).Value = 1.685235086126111
$ws.Cells.Item(2, 18).Value = 15.167115775135
$ws.Cells.Item(2, 19).Value = 0.04532758323745775
$ws.Cells.Item(2, 20).Value = 0.04532758323745775

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Leap2"
$ws.Cells.Item(3, 3).Value = "Ghsr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.093697666666667
$ws.Cells.Item(3, 8).Value = 9.281093
$ws.Cells.Item(3, 9).Value = 0.3054078832615814
$ws.Cells.Item(3, 10).Value = 0.3054078832615814
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.371854333333333
$ws.Cells.Item(3, 14).Value = 7.115563
$ws.Cells.Item(3, 15).Value = 0.6462308875194944
$ws.Cells.Item(3, 16).Value = 0.6462308875194943
$ws.Cells.Item(3, 17).Value = 7.337800216706556
$ws.Cells.Item(3, 18).Value = 66.040201950359
$ws.Cells.Item(3, 19).Value = 0.1973640074555819
$ws.Cells.Item(3, 20).Value = 0.1973640074555818

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Leap2"
$ws.Cells.Item(4, 3).Value = "Ghsr"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.093697666666667
$ws.Cells.Item(4, 8).Value = 9.281093
$ws.Cells.Item(4, 9).Value = 0.3054078832615814
$ws.Cells.Item(4, 10).Value = 0.3054078832615814
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.7537033333333333
$ws.Cells.Item(4, 14).Value = 2.26111
$ws.Cells.Item(4, 15).Value = 0.205352566210039
$ws.Cells.Item(4, 16).Value = 0.205352566210039
$ws.Cells.Item(4, 17).Value = 2.331730243692222
$ws.Cells.Item(4, 18).Value = 20.98557219323
$ws.Cells.Item(4, 19).Value = 0.06271629256854176
$ws.Cells.Item(4, 20).Value = 0.06271629256854176

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Leap2"
$ws.Cells.Item(5, 3).Value = "Ghsr"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.670935
$ws.Cells.Item(5, 8).Value = 8.012805
$ws.Cells.Item(5, 9).Value = 0.2636730193348796
$ws.Cells.Item(5, 10).Value = 0.2636730193348796
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.5447316666666667
$ws.Cells.Item(5, 14).Value = 1.634195
$ws.Cells.Item(5, 15).Value = 0.1484165462704666
$ws.Cells.Item(5, 16).Value = 0.1484165462704666
$ws.Cells.Item(5, 17).Value = 1.454942874108333
$ws.Cells.Item(5, 18).Value = 13.094485866975
$ws.Cells.Item(5, 19).Value = 0.03913343887438879
$ws.Cells.Item(5, 20).Value = 0.03913343887438879

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Leap2"
$ws.Cells.Item(6, 3).Value = "Ghsr"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.670935
$ws.Cells.Item(6, 8).Value = 8.012805
$ws.Cells.Item(6, 9).Value = 0.2636730193348796
$ws.Cells.Item(6, 10).Value = 0.2636730193348796
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.371854333333333
$ws.Cells.Item(6, 14).Value = 7.115563
$ws.Cells.Item(6, 15).Value = 0.6462308875194944
$ws.Cells.Item(6, 16).Value = 0.6462308875194943
$ws.Cells.Item(6, 17).Value = 6.335068753801667
$ws.Cells.Item(6, 18).Value = 57.015618784215
$ws.Cells.Item(6, 19).Value = 0.170393649299724
$ws.Cells.Item(6, 20).Value = 0.170393649299724

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Leap2"
$ws.Cells.Item(7, 3).Value = "Ghsr"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.670935
$ws.Cells.Item(7, 8).Value = 8.012805
$ws.Cells.Item(7, 9).Value = 0.2636730193348796
$ws.Cells.Item(7, 10).Value = 0.2636730193348796
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.7537033333333333
$ws.Cells.Item(7, 14).Value = 2.26111
$ws.Cells.Item(7, 15).Value = 0.205352566210039
$ws.Cells.Item(7, 16).Value = 0.205352566210039
$ws.Cells.Item(7, 17).Value = 2.013092612616667
$ws.Cells.Item(7, 18).Value = 18.11783351355
$ws.Cells.Item(7, 19).Value = 0.05414593116076676
$ws.Cells.Item(7, 20).Value = 0.05414593116076676

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Leap2"
$ws.Cells.Item(8, 3).Value = "Ghsr"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.819128666666666
$ws.Cells.Item(8, 8).Value = 8.457386
$ws.Cells.Item(8, 9).Value = 0.27830260468095
$ws.Cells.Item(8, 10).Value = 0.27830260468095
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.5447316666666667
$ws.Cells.Item(8, 14).Value = 1.634195
$ws.Cells.Item(8, 15).Value = 0.1484165462704666
$ws.Cells.Item(8, 16).Value = 0.1484165462704666
$ws.Cells.Item(8, 17).Value = 1.535668657141111
$ws.Cells.Item(8, 18).Value = 13.82101791427
$ws.Cells.Item(8, 19).Value = 0.04130471140482159
$ws.Cells.Item(8, 20).Value = 0.04130471140482159

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Leap2"
$ws.Cells.Item(9, 3).Value = "Ghsr"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.819128666666666
$ws.Cells.Item(9, 8).Value = 8.457386
$ws.Cells.Item(9, 9).Value = 0.27830260468095
$ws.Cells.Item(9, 10).Value = 0.27830260468095
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.371854333333333
$ws.Cells.Item(9, 14).Value = 7.115563
$ws.Cells.Item(9, 15).Value = 0.6462308875194944
$ws.Cells.Item(9, 16).Value = 0.6462308875194943
$ws.Cells.Item(9, 17).Value = 6.686562544257555
$ws.Cells.Item(9, 18).Value = 60.17906289831799
$ws.Cells.Item(9, 19).Value = 0.1798477392219573
$ws.Cells.Item(9, 20).Value = 0.1798477392219573

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Leap2"
$ws.Cells.Item(10, 3).Value = "Ghsr"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 2.819128666666666
$ws.Cells.Item(10, 8).Value = 8.457386
$ws.Cells.Item(10, 9).Value = 0.27830260468095
$ws.Cells.Item(10, 10).Value = 0.27830260468095
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.7537033333333333
$ws.Cells.Item(10, 14).Value = 2.26111
$ws.Cells.Item(10, 15).Value = 0.205352566210039
$ws.Cells.Item(10, 16).Value = 0.205352566210039
$ws.Cells.Item(10, 17).Value = 2.124786673162222
$ws.Cells.Item(10, 18).Value = 19.12308005846
$ws.Cells.Item(10, 19).Value = 0.0571501540541711
$ws.Cells.Item(10, 20).Value = 0.05715015405417111

# Row 11
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Leap2"
$ws.Cells.Item(11, 3).Value = "Ghsr"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.545963
$ws.Cells.Item(11, 8).Value = 4.637888999999999
$ws.Cells.Item(11, 9).Value = 0.152616492722589
$ws.Cells.Item(11, 10).Value = 0.152616492722589
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.5447316666666667
$ws.Cells.Item(11, 14).Value = 1.634195
$ws.Cells.Item(11, 15).Value = 0.1484165462704666
$ws.Cells.Item(11, 16).Value = 0.1484165462704666
$ws.Cells.Item(11, 17).Value = 0.842135001595
$ws.Cells.Item(11, 18).Value = 7.579215014354999
$ws.Cells.Item(11, 19).Value = 0.02265081275379846
$ws.Cells.Item(11, 20).Value = 0.02265081275379846

# Row 12
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value = "Leap2"
$ws.Cells.Item(12, 3).Value = "Ghsr"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.545963
$ws.Cells.Item(12, 8).Value = 4.637888999999999
$ws.Cells.Item(12, 9).Value = 0.152616492722589
$ws.Cells.Item(12, 10).Value = 0.152616492722589
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 2.371854333333333
$ws.Cells.Item(12, 14).Value = 7.115563
$ws.Cells.Item(12, 15).Value = 0.6462308875194944
$ws.Cells.Item(12, 16).Value = 0.6462308875194943
$ws.Cells.Item(12, 17).Value = 3.666799040722999
$ws.Cells.Item(12, 18).Value = 33.00119136650699
$ws.Cells.Item(12, 19).Value = 0.09862549154223117
$ws.Cells.Item(12, 20).Value = 0.09862549154223116

# Row 13
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value = "Leap2"
$ws.Cells.Item(13, 3).Value = "Ghsr"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1.545963
$ws.Cells.Item(13, 8).Value = 4.637888999999999
$ws.Cells.Item(13, 9).Value = 0.152616492722589
$ws.Cells.Item(13, 10).Value = 0.152616492722589
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.7537033333333333
$ws.Cells.Item(13, 14).Value = 2.26111
$ws.Cells.Item(13, 15).Value = 0.205352566210039
$ws.Cells.Item(13, 16).Value = 0.205352566210039
$ws.Cells.Item(13, 17).Value = 1.16519746631
$ws.Cells.Item(13, 18).Value = 10.48677719679
$ws.Cells.Item(13, 19).Value = 0.0313401884265594
$ws.Cells.Item(13, 20).Value = 0.0313401884265594

# Remove the now-obsolete rows 14-17 (Resolving-Mac -> * target rows dropped)
$ws.Range("A14:T17").Delete()
